$wb = $excel.ActiveWorkbook

# This script updates Kujata profit/price data cells (columns H-N) across
# several sheets, reflecting refreshed market data from the scheduled runner.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1983.7894
$ws.Range("I40").Value = 2766.6667
$ws.Range("J40").Value = 1622.4615
$ws.Range("K40").Value = 2766.6667
$ws.Range("L40").Value = 1622.4615
$ws.Range("M40").Value = -2591.6667
$ws.Range("N40").Value = -1972.4615

$ws.Range("H69").Value = 3980
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3980
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11940
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -13688

$ws.Range("H72").Value = 3980
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3980
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 35820
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -44556

$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 4000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5872

$ws.Range("H76").Value = 5848.3687
$ws.Range("I76").Value = 5087.8
$ws.Range("J76").Value = 6120
$ws.Range("K76").Value = 5087.8
$ws.Range("L76").Value = 6120
$ws.Range("M76").Value = -4772.8
$ws.Range("N76").Value = -6750

$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 20000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -29360

$ws.Range("H79").Value = 5848.3687
$ws.Range("I79").Value = 5087.8
$ws.Range("J79").Value = 6120
$ws.Range("K79").Value = 5087.8
$ws.Range("L79").Value = 6120
$ws.Range("M79").Value = -3995.8
$ws.Range("N79").Value = -8304

$ws.Range("H80").Value = 367.5
$ws.Range("I80").Value = 149.6
$ws.Range("J80").Value = 523.1429000000001
$ws.Range("K80").Value = 448.8
$ws.Range("L80").Value = 1569.4287
$ws.Range("M80").Value = 549.2
$ws.Range("N80").Value = -3565.4287

$ws.Range("H83").Value = 367.5
$ws.Range("I83").Value = 149.6
$ws.Range("J83").Value = 523.1429000000001
$ws.Range("K83").Value = 1346.4
$ws.Range("L83").Value = 4708.2861
$ws.Range("M83").Value = 3645.6
$ws.Range("N83").Value = -14692.2861

$ws.Range("H88").Value = 1237067.8
$ws.Range("I88").Value = 603
$ws.Range("J88").Value = 1374452.8
$ws.Range("K88").Value = 603
$ws.Range("L88").Value = 1374452.8
$ws.Range("M88").Value = -197
$ws.Range("N88").Value = -1375264.8

$ws.Range("H91").Value = 1237067.8
$ws.Range("I91").Value = 603
$ws.Range("J91").Value = 1374452.8
$ws.Range("K91").Value = 603
$ws.Range("L91").Value = 1374452.8
$ws.Range("M91").Value = 801
$ws.Range("N91").Value = -1377260.8

$ws.Range("H137").Value = 1185.1538
$ws.Range("I137").Value = 941.1053000000001
$ws.Range("J137").Value = 1847.5714
$ws.Range("K137").Value = 2823.3159
$ws.Range("L137").Value = 5542.7142
$ws.Range("M137").Value = -273.3159000000001
$ws.Range("N137").Value = -10642.7142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 9261863
$ws.Range("I102").Value = 10419227
$ws.Range("K102").Value = 10419227
$ws.Range("M102").Value = -10417605

$ws.Range("H124").Value = 3482.5
$ws.Range("J124").Value = 3482.5
$ws.Range("L124").Value = 3482.5
$ws.Range("N124").Value = -13302.5

$ws.Range("H135").Value = 21000
$ws.Range("J135").Value = 21000
$ws.Range("L135").Value = 21000
$ws.Range("N135").Value = -31140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1741
$ws.Range("I31").Value = 1412.875
$ws.Range("K31").Value = 1412.875
$ws.Range("M31").Value = -1117.875

$ws.Range("H34").Value = 1741
$ws.Range("I34").Value = 1412.875
$ws.Range("K34").Value = 1412.875
$ws.Range("M34").Value = -1210.875

$ws.Range("H134").Value = 1753.8125
$ws.Range("I134").Value = 1866
$ws.Range("J134").Value = 1467.1111
$ws.Range("K134").Value = 5598
$ws.Range("L134").Value = 4401.3333
$ws.Range("M134").Value = -3063
$ws.Range("N134").Value = -9471.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1327.8148
$ws.Range("I5").Value = 1406.3334
$ws.Range("J5").Value = 699.6667
$ws.Range("K5").Value = 4219.0002
$ws.Range("L5").Value = 2099.0001
$ws.Range("M5").Value = -4107.0002
$ws.Range("N5").Value = -2323.0001

$ws.Range("H131").Value = 16950462
$ws.Range("J131").Value = 1407.0555
$ws.Range("L131").Value = 4221.166499999999
$ws.Range("N131").Value = -14301.1665

$ws.Range("H135").Value = 1327.8148
$ws.Range("I135").Value = 1406.3334
$ws.Range("J135").Value = 699.6667
$ws.Range("K135").Value = 12657.0006
$ws.Range("L135").Value = 6297.0003
$ws.Range("M135").Value = -10122.0006
$ws.Range("N135").Value = -11367.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6281.6665
$ws.Range("I80").Value = 6790
$ws.Range("J80").Value = 6180
$ws.Range("K80").Value = 6790
$ws.Range("L80").Value = 6180
$ws.Range("M80").Value = -5792
$ws.Range("N80").Value = -8176

$ws.Range("H83").Value = 6281.6665
$ws.Range("I83").Value = 6790
$ws.Range("J83").Value = 6180
$ws.Range("K83").Value = 33950
$ws.Range("L83").Value = 30900
$ws.Range("M83").Value = -28958
$ws.Range("N83").Value = -40884

$ws.Range("H126").Value = 1636.3462
$ws.Range("I126").Value = 1466.5294
$ws.Range("J126").Value = 1957.1111
$ws.Range("K126").Value = 4399.5882
$ws.Range("L126").Value = 5871.3333
$ws.Range("M126").Value = -1929.5882
$ws.Range("N126").Value = -10811.3333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1571.4286
$ws.Range("I22").Value = 499
$ws.Range("K22").Value = 499
$ws.Range("M22").Value = -204

$ws.Range("H27").Value = 1571.4286
$ws.Range("I27").Value = 499
$ws.Range("K27").Value = 499
$ws.Range("M27").Value = -392

$ws.Range("H93").Value = 886.0294
$ws.Range("I93").Value = 604.5
$ws.Range("J93").Value = 2199.8333
$ws.Range("K93").Value = 604.5
$ws.Range("L93").Value = 2199.8333
$ws.Range("M93").Value = 643.5
$ws.Range("N93").Value = -4695.8333

$ws.Range("H96").Value = 14000
$ws.Range("J96").Value = 14000
$ws.Range("L96").Value = 14000
$ws.Range("N96").Value = -19492

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12133
$ws.Range("J52").Value = 12133
$ws.Range("L52").Value = 12133
$ws.Range("N52").Value = -12585

$ws.Range("H81").Value = 650
$ws.Range("I81").Value = 650
$ws.Range("K81").Value = 1300
$ws.Range("M81").Value = -239

$ws.Range("H84").Value = 650
$ws.Range("I84").Value = 650
$ws.Range("K84").Value = 6500
$ws.Range("M84").Value = -1196
